$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-1.0049686186707034, 1.2403050310127643, 0.3427709853979849, 3.0544146098810216)"
$ws.Range("C2").Value = "JSU(-1.4574060114848337, 1.1305676093944697, 2.054772917862269, 4.245666484153441)"
$ws.Range("D2").Value = "JSU(-0.8394915002860953, 1.0499210393966387, 0.637609969180015, 2.598650886635727)"
$ws.Range("E2").Value = "NIG(1.5981080917869037, 1.2579251679529966, 3.0970105522520885, 6.070250695725709)"
